# Third tab ("Лист1" / Sheet1) of the audiences workbook:
#  - widen column A slightly (custom column width)
#  - re-tag two existing "any" audience rows (B11, B16) as "Компьютерный класс"
#  - append a brand-new audience row (37): "Спортивный зал" / "Спортивные снаряды"
#  - leave the selection on B16, the cell that was last edited

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width -------------------------------------------------------
# Target OOXML col width is 15.109375 characters. This runtime's ColumnWidth
# setter quantizes to 1/6-character steps, so that exact figure can't be
# reproduced bit-for-bit; 14.333333 is the input that lands on the closest
# reachable stored width (15.1666...).
$ws.Columns.Item(1).ColumnWidth = 14.333333

# --- New row 37 first, so its strings land before the "Компьютерный класс" one ---
$ws.Range("A37").Value = "Спортивный зал"
$ws.Range("B37").Value = "Спортивные снаряды"

# --- Re-classify rows 11 and 16 from "any" to "Компьютерный класс" --------
$ws.Range("B11").Value = "Компьютерный класс"
$ws.Range("B16").Value = "Компьютерный класс"

# --- Match the final on-screen selection -----------------------------------
$ws.Range("B16").Select()
